# Updates cryptos list cell values to match the latest scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.333.22'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '2.498.89'
$ws.Range('E3').Value = '  -2.57%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('E7').Value = '  +1.76%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.535'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.36%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.50'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('E13').Value = '  -4.21%  '
$ws.Range('D14').Value = '2.889.57'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').Value = '2.531.84'
$ws.Range('E15').Value = '  -2.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.852'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.65%  '
$ws.Range('D18').Value = '42.340.17'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('E20').Value = '  -1.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.73%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  +10.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.15'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.25%  '
$ws.Range('E34').Value = '  -1.54%  '
$ws.Range('E35').Value = '  -1.94%  '
$ws.Range('E36').Value = '  -3.88%  '
$ws.Range('E37').Value = '  -5.65%  '
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.86'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.89'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.34%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0300'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.029.17'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '84.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.74%  '
$ws.Range('E48').Value = '  -2.61%  '
$ws.Range('D49').Value = '2.746.66'
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.95%  '
$ws.Range('E51').Value = '  -0.08%  '
